$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: split "JOHP; " into proofErr-wrapped "JOHP;" + separate " " run ---
$p3 = $d.Paragraphs(3)
if ($p3.Range.Text -ne "JOHP; `r") {
    throw "unexpected paragraph 3 text: [$($p3.Range.Text)]"
}
$xmlJohp = '<w:p ' + $wNs + '><w:proofErr w:type="gramStart"/><w:r><w:t>JOHP;</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$p3.Range.InsertXML($xmlJohp)

# --- Step 2: split "work and stress;" into "work and " + proofErr-wrapped "stress;" ---
#             and simultaneously insert the new "anxiety stress and coping" paragraph
#             right after it (this paragraph moved up from later in the document).
$p4 = $d.Paragraphs(4)
if ($p4.Range.Text -ne "work and stress;`r") {
    throw "unexpected paragraph 4 text: [$($p4.Range.Text)]"
}
$xmlWork = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">work and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>stress;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>' + `
           '<w:p ' + $wNs + '><w:r><w:t>anxiety stress and coping</w:t></w:r></w:p>'
$p4.Range.InsertXML($xmlWork)

# --- Step 3: remove the old "anxiety stress and coping" paragraph that used to sit
#             right after "Potential Options" ---
$removed = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "anxiety stress and coping`r" -and $p.Previous().Range.Text -eq "Potential Options`r") {
        $p.Range.Delete()
        $removed = $true
        break
    }
}
if (-not $removed) {
    throw "could not find the old 'anxiety stress and coping' paragraph to remove"
}
